$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 167.76923
$ws.Range("I9").Value = 153.44444
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 153.44444
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 15.55556000000001
$ws.Range("N9").Value = -538

$ws.Range("H129").Value = 2779.5894
$ws.Range("I129").Value = 6075.222
$ws.Range("J129").Value = 1218.5
$ws.Range("K129").Value = 18225.666
$ws.Range("L129").Value = 3655.5
$ws.Range("M129").Value = -13225.666
$ws.Range("N129").Value = -13655.5

$ws.Range("H132").Value = 3574844
$ws.Range("I132").Value = 3734870.2
$ws.Range("K132").Value = 11204610.6
$ws.Range("M132").Value = -11202080.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 16.666666
$ws.Range("I25").Value = 16.666666
$ws.Range("K25").Value = 16.666666
$ws.Range("M25").Value = 385.333334

$ws.Range("H32").Value = 21260.035
$ws.Range("I32").Value = 3958.1357
$ws.Range("K32").Value = 3958.1357
$ws.Range("M32").Value = -3671.1357

$ws.Range("H45").Value = 1268.3125
$ws.Range("I45").Value = 854.875
$ws.Range("J45").Value = 1681.75
$ws.Range("K45").Value = 854.875
$ws.Range("L45").Value = 1681.75
$ws.Range("M45").Value = -477.875
$ws.Range("N45").Value = -2435.75

$ws.Range("H61").Value = 1682.34
$ws.Range("I61").Value = 1045.3572
$ws.Range("J61").Value = 2493.0454
$ws.Range("K61").Value = 1045.3572
$ws.Range("L61").Value = 2493.0454
$ws.Range("M61").Value = -833.3571999999999
$ws.Range("N61").Value = -2917.0454

$ws.Range("H74").Value = 831.4286
$ws.Range("I74").Value = 807.4
$ws.Range("J74").Value = 844.7778
$ws.Range("K74").Value = 807.4
$ws.Range("L74").Value = 844.7778
$ws.Range("M74").Value = 66.60000000000002
$ws.Range("N74").Value = -2592.7778

$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36748

$ws.Range("H77").Value = 831.4286
$ws.Range("I77").Value = 807.4
$ws.Range("J77").Value = 844.7778
$ws.Range("K77").Value = 4037
$ws.Range("L77").Value = 4223.889
$ws.Range("M77").Value = 331
$ws.Range("N77").Value = -12959.889

$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -113736

$ws.Range("H110").Value = 71500730
$ws.Range("I110").Value = 83417330
$ws.Range("J110").Value = 1156.5
$ws.Range("K110").Value = 83417330
$ws.Range("L110").Value = 1156.5
$ws.Range("M110").Value = -83415285
$ws.Range("N110").Value = -5246.5

$ws.Range("H124").Value = 21803.625
$ws.Range("J124").Value = 21803.625
$ws.Range("L124").Value = 21803.625
$ws.Range("N124").Value = -31623.625

$ws.Range("H136").Value = 1682.34
$ws.Range("I136").Value = 1045.3572
$ws.Range("J136").Value = 2493.0454
$ws.Range("K136").Value = 3136.0716
$ws.Range("L136").Value = 7479.1362
$ws.Range("M136").Value = -586.0715999999998
$ws.Range("N136").Value = -12579.1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 23169.826
$ws.Range("I20").Value = 33865.902
$ws.Range("J20").Value = 1064.6
$ws.Range("K20").Value = 33865.902
$ws.Range("L20").Value = 1064.6
$ws.Range("M20").Value = -33618.902
$ws.Range("N20").Value = -1558.6

$ws.Range("H99").Value = 2530
$ws.Range("I99").Value = 2790
$ws.Range("J99").Value = 2478
$ws.Range("K99").Value = 2790
$ws.Range("L99").Value = 2478
$ws.Range("M99").Value = -1292
$ws.Range("N99").Value = -5474

$ws.Range("H102").Value = 6000
$ws.Range("I102").Value = 6000
$ws.Range("K102").Value = 6000
$ws.Range("M102").Value = -2755

$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 15900
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 29800
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 29800
$ws.Range("M39").Value = -1609
$ws.Range("N39").Value = -30582

$ws.Range("H49").Value = 15900
$ws.Range("I49").Value = 2000
$ws.Range("J49").Value = 29800
$ws.Range("K49").Value = 2000
$ws.Range("L49").Value = 29800
$ws.Range("M49").Value = -1818
$ws.Range("N49").Value = -30164

$ws.Range("H58").Value = 1084.3934
$ws.Range("I58").Value = 956.1818
$ws.Range("J58").Value = 2259.6667
$ws.Range("K58").Value = 956.1818
$ws.Range("L58").Value = 2259.6667
$ws.Range("M58").Value = -753.1818
$ws.Range("N58").Value = -2665.6667

$ws.Range("H62").Value = 2730
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -2126
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 2730
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -10630
$ws.Range("N65").Value = -19740

$ws.Range("H68").Value = 19328
$ws.Range("J68").Value = 19328
$ws.Range("L68").Value = 19328
$ws.Range("N68").Value = -20826

$ws.Range("H71").Value = 19328
$ws.Range("J71").Value = 19328
$ws.Range("L71").Value = 57984
$ws.Range("N71").Value = -65472

$ws.Range("H74").Value = 20224.857
$ws.Range("J74").Value = 26814.8
$ws.Range("L74").Value = 26814.8
$ws.Range("N74").Value = -28562.8

$ws.Range("H77").Value = 20224.857
$ws.Range("J77").Value = 26814.8
$ws.Range("L77").Value = 80444.39999999999
$ws.Range("N77").Value = -89180.39999999999

$ws.Range("H132").Value = 2294.6826
$ws.Range("I132").Value = 2226.8333
$ws.Range("J132").Value = 2511.8
$ws.Range("K132").Value = 6680.499899999999
$ws.Range("L132").Value = 7535.400000000001
$ws.Range("M132").Value = -4150.499899999999
$ws.Range("N132").Value = -12595.4

$ws.Range("H134").Value = 1103.0834
$ws.Range("I134").Value = 1043.7
$ws.Range("K134").Value = 3131.1
$ws.Range("M134").Value = -596.1000000000004

$ws.Range("H136").Value = 1084.3934
$ws.Range("I136").Value = 956.1818
$ws.Range("J136").Value = 2259.6667
$ws.Range("K136").Value = 2868.5454
$ws.Range("L136").Value = 6779.000100000001
$ws.Range("M136").Value = -318.5454
$ws.Range("N136").Value = -11879.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 450

$ws.Range("H131").Value = 928.6559
$ws.Range("I131").Value = 639.8182
$ws.Range("J131").Value = 967.40247
$ws.Range("K131").Value = 1919.4546
$ws.Range("L131").Value = 2902.20741
$ws.Range("M131").Value = 3120.5454
$ws.Range("N131").Value = -12982.20741

$ws.Range("H141").Value = 2242.7222
$ws.Range("I141").Value = 1274.5385
$ws.Range("J141").Value = 4760
$ws.Range("K141").Value = 3823.6155
$ws.Range("L141").Value = 14280
$ws.Range("M141").Value = 1356.3845
$ws.Range("N141").Value = -24640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 13133.333
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 13133.333
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 13133.333
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -13651.333

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900

$ws.Range("H126").Value = 14708651
$ws.Range("I126").Value = 3691.6667
$ws.Range("J126").Value = 58823530
$ws.Range("K126").Value = 11075.0001
$ws.Range("L126").Value = 176470590
$ws.Range("M126").Value = -8605.000100000001
$ws.Range("N126").Value = -176475530

$ws.Range("H131").Value = 40663.332
$ws.Range("J131").Value = 40663.332
$ws.Range("L131").Value = 40663.332
$ws.Range("N131").Value = -50743.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 14000
$ws.Range("J38").Value = 14000
$ws.Range("L38").Value = 14000
$ws.Range("N38").Value = -14820

$ws.Range("H45").Value = 6036.727
$ws.Range("I45").Value = 4166
$ws.Range("J45").Value = 6909.7334
$ws.Range("K45").Value = 4166
$ws.Range("L45").Value = 6909.7334
$ws.Range("M45").Value = -3759
$ws.Range("N45").Value = -7723.7334

$ws.Range("H48").Value = 15041
$ws.Range("I48").Value = 15041
$ws.Range("K48").Value = 15041
$ws.Range("M48").Value = -14380

$ws.Range("H122").Value = 2532.2
$ws.Range("I122").Value = 2600.4
$ws.Range("J122").Value = 2395.8
$ws.Range("K122").Value = 7801.200000000001
$ws.Range("L122").Value = 7187.400000000001
$ws.Range("M122").Value = -5351.200000000001
$ws.Range("N122").Value = -12087.4

$ws.Range("H127").Value = 28403.75
$ws.Range("J127").Value = 28403.75
$ws.Range("L127").Value = 28403.75
$ws.Range("N127").Value = -38323.75

$ws.Range("H136").Value = 934.9787
$ws.Range("I136").Value = 896.4545000000001
$ws.Range("K136").Value = 2689.3635
$ws.Range("M136").Value = -139.3635000000004
